$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 478.8
$ws.Cells.Item(28, 9).Value = 298.33334
$ws.Cells.Item(28, 10).Value = 749.5
$ws.Cells.Item(28, 11).Value = 298.33334
$ws.Cells.Item(28, 12).Value = 749.5
$ws.Cells.Item(28, 13).Value = 186.66666
$ws.Cells.Item(28, 14).Value = -1719.5
$ws.Cells.Item(51, 8).Value = 3596
$ws.Cells.Item(51, 10).Value = 3596
$ws.Cells.Item(51, 12).Value = 3596
$ws.Cells.Item(51, 14).Value = -4564
$ws.Cells.Item(86, 8).Value = 411461
$ws.Cells.Item(86, 10).Value = 1125
$ws.Cells.Item(86, 12).Value = 1125
$ws.Cells.Item(86, 14).Value = -3371
$ws.Cells.Item(89, 8).Value = 411461
$ws.Cells.Item(89, 10).Value = 1125
$ws.Cells.Item(89, 12).Value = 5625
$ws.Cells.Item(89, 14).Value = -16857
$ws.Cells.Item(137, 8).Value = 1586.591
$ws.Cells.Item(137, 10).Value = 3283.3333
$ws.Cells.Item(137, 12).Value = 9849.999899999999
$ws.Cells.Item(137, 14).Value = -14949.9999
$ws.Cells.Item(138, 8).Value = 4064.6572
$ws.Cells.Item(138, 10).Value = 3059.9644
$ws.Cells.Item(138, 12).Value = 9179.893199999999
$ws.Cells.Item(138, 14).Value = -19459.8932

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4754.375
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 13).ClearContents()
$ws.Cells.Item(61, 8).Value = 2250.3333
$ws.Cells.Item(61, 9).Value = 1459.0278
$ws.Cells.Item(61, 10).Value = 6998.1665
$ws.Cells.Item(61, 11).Value = 1459.0278
$ws.Cells.Item(61, 12).Value = 6998.1665
$ws.Cells.Item(61, 13).Value = -1247.0278
$ws.Cells.Item(61, 14).Value = -7422.1665
$ws.Cells.Item(74, 8).Value = 1827.9286
$ws.Cells.Item(74, 9).Value = 1592.8572
$ws.Cells.Item(74, 10).Value = 2063
$ws.Cells.Item(74, 11).Value = 1592.8572
$ws.Cells.Item(74, 12).Value = 2063
$ws.Cells.Item(74, 13).Value = -718.8571999999999
$ws.Cells.Item(74, 14).Value = -3811
$ws.Cells.Item(77, 8).Value = 1827.9286
$ws.Cells.Item(77, 9).Value = 1592.8572
$ws.Cells.Item(77, 10).Value = 2063
$ws.Cells.Item(77, 11).Value = 7964.286
$ws.Cells.Item(77, 12).Value = 10315
$ws.Cells.Item(77, 13).Value = -3596.286
$ws.Cells.Item(77, 14).Value = -19051
$ws.Cells.Item(97, 8).Value = 1164.5
$ws.Cells.Item(97, 9).Value = 1164.5
$ws.Cells.Item(97, 11).Value = 1164.5
$ws.Cells.Item(97, 13).Value = -668.5
$ws.Cells.Item(109, 8).Value = 85000
$ws.Cells.Item(109, 10).Value = 85000
$ws.Cells.Item(109, 12).Value = 85000
$ws.Cells.Item(109, 14).Value = -87774
$ws.Cells.Item(114, 8).Value = 9749
$ws.Cells.Item(114, 10).Value = 9749
$ws.Cells.Item(114, 12).Value = 9749
$ws.Cells.Item(114, 14).Value = -18427
$ws.Cells.Item(132, 8).Value = 1460.0233
$ws.Cells.Item(132, 9).Value = 890.8333
$ws.Cells.Item(132, 11).Value = 2672.4999
$ws.Cells.Item(132, 13).Value = -142.4998999999998
$ws.Cells.Item(133, 8).Value = 35000
$ws.Cells.Item(133, 10).Value = 35000
$ws.Cells.Item(133, 12).Value = 35000
$ws.Cells.Item(133, 14).Value = -40060
$ws.Cells.Item(136, 8).Value = 2250.3333
$ws.Cells.Item(136, 9).Value = 1459.0278
$ws.Cells.Item(136, 10).Value = 6998.1665
$ws.Cells.Item(136, 11).Value = 4377.0834
$ws.Cells.Item(136, 12).Value = 20994.4995
$ws.Cells.Item(136, 13).Value = -1827.0834
$ws.Cells.Item(136, 14).Value = -26094.4995

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 2044.6842
$ws.Cells.Item(107, 9).Value = 1842.6666
$ws.Cells.Item(107, 11).Value = 1842.6666
$ws.Cells.Item(107, 13).Value = 77.33339999999998
$ws.Cells.Item(134, 8).Value = 4682.2446
$ws.Cells.Item(134, 9).Value = 5199.6855
$ws.Cells.Item(134, 11).Value = 15599.0565
$ws.Cells.Item(134, 13).Value = -13064.0565

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1782.5758
$ws.Cells.Item(31, 9).Value = 1780.1154
$ws.Cells.Item(31, 10).Value = 1791.7142
$ws.Cells.Item(31, 11).Value = 1780.1154
$ws.Cells.Item(31, 12).Value = 1791.7142
$ws.Cells.Item(31, 13).Value = -1485.1154
$ws.Cells.Item(31, 14).Value = -2381.7142
$ws.Cells.Item(34, 8).Value = 1782.5758
$ws.Cells.Item(34, 9).Value = 1780.1154
$ws.Cells.Item(34, 10).Value = 1791.7142
$ws.Cells.Item(34, 11).Value = 1780.1154
$ws.Cells.Item(34, 12).Value = 1791.7142
$ws.Cells.Item(34, 13).Value = -1578.1154
$ws.Cells.Item(34, 14).Value = -2195.7142
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 12).ClearContents()
$ws.Cells.Item(48, 14).Value = 0
$ws.Cells.Item(105, 8).Value = 1660.6875
$ws.Cells.Item(105, 9).Value = 1348
$ws.Cells.Item(105, 10).Value = 2348.6
$ws.Cells.Item(105, 11).Value = 1348
$ws.Cells.Item(105, 12).Value = 2348.6
$ws.Cells.Item(105, 13).Value = 399
$ws.Cells.Item(105, 14).Value = -5842.6
$ws.Cells.Item(134, 8).Value = 966.569
$ws.Cells.Item(134, 9).Value = 885.8298
$ws.Cells.Item(134, 10).Value = 1311.5454
$ws.Cells.Item(134, 11).Value = 2657.4894
$ws.Cells.Item(134, 12).Value = 3934.6362
$ws.Cells.Item(134, 13).Value = -122.4893999999999
$ws.Cells.Item(134, 14).Value = -9004.636200000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 459.6
$ws.Cells.Item(18, 9).Value = 432.8889
$ws.Cells.Item(18, 10).Value = 700
$ws.Cells.Item(18, 11).Value = 1298.6667
$ws.Cells.Item(18, 12).Value = 2100
$ws.Cells.Item(18, 13).Value = -1129.6667
$ws.Cells.Item(18, 14).Value = -2438
$ws.Cells.Item(38, 8).Value = 357
$ws.Cells.Item(38, 9).Value = 99.25
$ws.Cells.Item(38, 10).Value = 700.6667
$ws.Cells.Item(38, 11).Value = 297.75
$ws.Cells.Item(38, 12).Value = 2102.0001
$ws.Cells.Item(38, 13).Value = 49.25
$ws.Cells.Item(38, 14).Value = -2796.0001
$ws.Cells.Item(68, 8).Value = 758.9091
$ws.Cells.Item(68, 9).Value = 759.2
$ws.Cells.Item(68, 10).Value = 758.6667
$ws.Cells.Item(68, 11).Value = 2277.6
$ws.Cells.Item(68, 12).Value = 2276.0001
$ws.Cells.Item(68, 13).Value = -1466.6
$ws.Cells.Item(68, 14).Value = -3898.0001
$ws.Cells.Item(69, 8).Value = 2500
$ws.Cells.Item(69, 9).Value = 2500
$ws.Cells.Item(69, 11).Value = 7500
$ws.Cells.Item(69, 13).Value = -6689
$ws.Cells.Item(71, 8).Value = 758.9091
$ws.Cells.Item(71, 9).Value = 759.2
$ws.Cells.Item(71, 10).Value = 758.6667
$ws.Cells.Item(71, 11).Value = 6832.8
$ws.Cells.Item(71, 12).Value = 6828.0003
$ws.Cells.Item(71, 13).Value = -2776.8
$ws.Cells.Item(71, 14).Value = -14940.0003
$ws.Cells.Item(72, 8).Value = 2500
$ws.Cells.Item(72, 9).Value = 2500
$ws.Cells.Item(72, 11).Value = 22500
$ws.Cells.Item(72, 13).Value = -18444
$ws.Cells.Item(74, 8).Value = 1800
$ws.Cells.Item(74, 9).Value = 1800
$ws.Cells.Item(74, 11).Value = 5400
$ws.Cells.Item(74, 13).Value = -4339
$ws.Cells.Item(77, 8).Value = 1800
$ws.Cells.Item(77, 9).Value = 1800
$ws.Cells.Item(77, 11).Value = 16200
$ws.Cells.Item(77, 13).Value = -10896
$ws.Cells.Item(131, 8).Value = 11153.985
$ws.Cells.Item(131, 10).Value = 12155.822
$ws.Cells.Item(131, 12).Value = 36467.466
$ws.Cells.Item(131, 14).Value = -46547.466
$ws.Cells.Item(132, 8).Value = 843.2857
$ws.Cells.Item(132, 9).Value = 634.6667
$ws.Cells.Item(132, 11).Value = 5712.0003
$ws.Cells.Item(132, 13).Value = -3182.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1035.3
$ws.Cells.Item(113, 10).Value = 1202.6
$ws.Cells.Item(113, 12).Value = 1202.6
$ws.Cells.Item(113, 14).Value = -5542.6
$ws.Cells.Item(132, 8).Value = 895995
$ws.Cells.Item(132, 9).Value = 1132363.5
$ws.Cells.Item(132, 10).Value = 3047.111
$ws.Cells.Item(132, 11).Value = 3397090.5
$ws.Cells.Item(132, 12).Value = 9141.332999999999
$ws.Cells.Item(132, 13).Value = -3394560.5
$ws.Cells.Item(132, 14).Value = -14201.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1863.2727
$ws.Cells.Item(61, 10).Value = 1810.625
$ws.Cells.Item(61, 12).Value = 1810.625
$ws.Cells.Item(61, 14).Value = -2214.625
$ws.Cells.Item(113, 8).Value = 1863.2727
$ws.Cells.Item(113, 10).Value = 1810.625
$ws.Cells.Item(113, 12).Value = 1810.625
$ws.Cells.Item(113, 14).Value = -6150.625
$ws.Cells.Item(132, 8).Value = 1591.1372
$ws.Cells.Item(132, 9).Value = 1087
$ws.Cells.Item(132, 10).Value = 2115.44
$ws.Cells.Item(132, 11).Value = 3261
$ws.Cells.Item(132, 12).Value = 6346.32
$ws.Cells.Item(132, 13).Value = -731
$ws.Cells.Item(132, 14).Value = -11406.32
$ws.Cells.Item(136, 8).Value = 2134.8865
$ws.Cells.Item(136, 9).Value = 1178.8286
$ws.Cells.Item(136, 11).Value = 3536.4858
$ws.Cells.Item(136, 13).Value = -986.4858000000004

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(75, 8).Value = 29999
$ws.Cells.Item(75, 10).Value = 29999
$ws.Cells.Item(75, 12).Value = 29999
$ws.Cells.Item(75, 14).Value = -31871
$ws.Cells.Item(78, 8).Value = 29999
$ws.Cells.Item(78, 10).Value = 29999
$ws.Cells.Item(78, 12).Value = 89997
$ws.Cells.Item(78, 14).Value = -99357
$ws.Cells.Item(132, 8).Value = 1711.9697
$ws.Cells.Item(132, 9).Value = 888.5217
$ws.Cells.Item(132, 10).Value = 3605.9
$ws.Cells.Item(132, 11).Value = 2665.5651
$ws.Cells.Item(132, 12).Value = 10817.7
$ws.Cells.Item(132, 13).Value = -135.5650999999998
$ws.Cells.Item(132, 14).Value = -15877.7
$ws.Cells.Item(136, 8).Value = 16341902
$ws.Cells.Item(136, 9).Value = 22224046
$ws.Cells.Item(136, 10).Value = 2611.6667
$ws.Cells.Item(136, 11).Value = 66672138
$ws.Cells.Item(136, 12).Value = 7835.000100000001
$ws.Cells.Item(136, 13).Value = -66669588
$ws.Cells.Item(136, 14).Value = -12935.0001
